# act tablas web jul25
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Data": add 2023 and 2022 rows at the top (pushing the
# existing years down) and extend the historical series back from
# 2005 through to 1985 at the bottom.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

# Insert two new rows right after the header row for 2023 / 2022.
$ws.Range("A2:A3").EntireRow.Insert()

# Column A holds the year as text (matches the rest of the column,
# which stores years as shared strings rather than numbers).
$yearCol = $ws.Range("A2:A40")
$yearCol.NumberFormat = "@"

$years = @(2023,2022,2021,2020,2019,2018,2017,2016,2015,2014,2013,2012,2011,2010,2009,2008,2007,2006,2005,2004,2003,2002,2001,2000,1999,1998,1997,1996,1995,1994,1993,1992,1991,1990,1989,1988,1987,1986,1985)
$values = @(18.5,18.4,17.7,17.6,18.6,18.4,18.1,18.1,17.8,18.4,18.4,18.8,18.6,18.2,19.4,19,17.4,15.8,15.6,15.7,15.1,13.6,13.1,12.7,12.3,14.9,11.9,11.6,10.1,10.5,12.4,11.9,12.9,14.6,15.9,15.7,16,15.9,14.5)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = "$($years[$i])"
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# ---------------------------------------------------------------
# Sheet "Metadata": add an "actualizacion" / "Julio 2025" row just
# before the "cita" row.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Metadata")
$ws2.Range("A9:B9").EntireRow.Insert()
$ws2.Range("A9").Value = "actualizacion"
$ws2.Range("B9").Value = "Julio 2025"

# The blank A1 cell is re-saved using the same single-space marker
# already used by B1.
$ws2.Range("A1").Value = " "
